# Cell updates derived from the crypto price/volume refresh diff.
# Each entry is (cell reference, new value, whether the text must be forced
# to remain a text value instead of being auto-parsed as a number by Excel).
$updates = @(
    @{ Cell = "D2"; Value = "41.536.86"; ForceText = $false }
    @{ Cell = "E2"; Value = "  +0.18%  "; ForceText = $false }
    @{ Cell = "D3"; Value = "2.489.80"; ForceText = $false }
    @{ Cell = "E3"; Value = "  +1.06%  "; ForceText = $false }
    @{ Cell = "D4"; Value = "0.996"; ForceText = $true }
    @{ Cell = "E4"; Value = "  -0.41%  "; ForceText = $false }
    @{ Cell = "D5"; Value = "312.51"; ForceText = $true }
    @{ Cell = "E5"; Value = "  +0.39%  "; ForceText = $false }
    @{ Cell = "D6"; Value = "94.19"; ForceText = $true }
    @{ Cell = "E6"; Value = "  +0.31%  "; ForceText = $false }
    @{ Cell = "D7"; Value = "0.545"; ForceText = $true }
    @{ Cell = "E7"; Value = "  -1.10%  "; ForceText = $false }
    @{ Cell = "D8"; Value = "0.998"; ForceText = $true }
    @{ Cell = "E8"; Value = "  -0.35%  "; ForceText = $false }
    @{ Cell = "D9"; Value = "0.497"; ForceText = $true }
    @{ Cell = "E9"; Value = "  -1.96%  "; ForceText = $false }
    @{ Cell = "D10"; Value = "32.66"; ForceText = $true }
    @{ Cell = "E10"; Value = "  -2.31%  "; ForceText = $false }
    @{ Cell = "E11"; Value = "  +0.26%  "; ForceText = $false }
    @{ Cell = "E12"; Value = "  +1.61%  "; ForceText = $false }
    @{ Cell = "D13"; Value = "2.872.74"; ForceText = $false }
    @{ Cell = "E13"; Value = "  +0.78%  "; ForceText = $false }
    @{ Cell = "D14"; Value = "6.86"; ForceText = $true }
    @{ Cell = "E14"; Value = "  -1.25%  "; ForceText = $false }
    @{ Cell = "D15"; Value = "15.45"; ForceText = $true }
    @{ Cell = "E15"; Value = "  +6.49%  "; ForceText = $false }
    @{ Cell = "D16"; Value = "2.480.15"; ForceText = $false }
    @{ Cell = "E16"; Value = "  -0.18%  "; ForceText = $false }
    @{ Cell = "D17"; Value = "0.762"; ForceText = $true }
    @{ Cell = "E17"; Value = "  -3.47%  "; ForceText = $false }
    @{ Cell = "D18"; Value = "41.599.69"; ForceText = $false }
    @{ Cell = "E18"; Value = "  +0.35%  "; ForceText = $false }
    @{ Cell = "D19"; Value = "6.32"; ForceText = $true }
    @{ Cell = "E19"; Value = "  -0.53%  "; ForceText = $false }
    @{ Cell = "D20"; Value = "0.0₃0921"; ForceText = $false }
    @{ Cell = "E20"; Value = "  +0.56%  "; ForceText = $false }
    @{ Cell = "D21"; Value = "70.89"; ForceText = $true }
    @{ Cell = "E21"; Value = "  +3.28%  "; ForceText = $false }
    @{ Cell = "D22"; Value = "11.15"; ForceText = $true }
    @{ Cell = "E22"; Value = "  -3.38%  "; ForceText = $false }
    @{ Cell = "D23"; Value = "235.80"; ForceText = $true }
    @{ Cell = "E23"; Value = "  -0.80%  "; ForceText = $false }
    @{ Cell = "E24"; Value = "  -2.49%  "; ForceText = $false }
    @{ Cell = "E25"; Value = "  +0.09%  "; ForceText = $false }
    @{ Cell = "D26"; Value = "1.90"; ForceText = $true }
    @{ Cell = "E26"; Value = "  -1.79%  "; ForceText = $false }
    @{ Cell = "D27"; Value = "24.54"; ForceText = $true }
    @{ Cell = "E27"; Value = "  -0.88%  "; ForceText = $false }
    @{ Cell = "E28"; Value = "  +1.42%  "; ForceText = $false }
    @{ Cell = "D29"; Value = "9.64"; ForceText = $true }
    @{ Cell = "E29"; Value = "  -0.86%  "; ForceText = $false }
    @{ Cell = "D30"; Value = "36.30"; ForceText = $true }
    @{ Cell = "E30"; Value = "  -0.09%  "; ForceText = $false }
    @{ Cell = "D31"; Value = "154.01"; ForceText = $true }
    @{ Cell = "E31"; Value = "  +0.54%  "; ForceText = $false }
    @{ Cell = "D32"; Value = "5.42"; ForceText = $true }
    @{ Cell = "E32"; Value = "  -3.82%  "; ForceText = $false }
    @{ Cell = "B33"; Value = "Celestia"; ForceText = $false }
    @{ Cell = "C33"; Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"; ForceText = $false }
    @{ Cell = "D33"; Value = "18.17"; ForceText = $true }
    @{ Cell = "E33"; Value = "  +6.20%  "; ForceText = $false }
    @{ Cell = "B34"; Value = "WEMIXToken"; ForceText = $false }
    @{ Cell = "C34"; Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"; ForceText = $false }
    @{ Cell = "D34"; Value = "2.56"; ForceText = $true }
    @{ Cell = "E34"; Value = "  -2.33%  "; ForceText = $false }
    @{ Cell = "B35"; Value = "Hedera"; ForceText = $false }
    @{ Cell = "C35"; Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"; ForceText = $false }
    @{ Cell = "D35"; Value = "0.0759"; ForceText = $true }
    @{ Cell = "E35"; Value = "  +0.93%  "; ForceText = $false }
    @{ Cell = "E36"; Value = "  -1.52%  "; ForceText = $false }
    @{ Cell = "D37"; Value = "2.98"; ForceText = $true }
    @{ Cell = "E37"; Value = "  -1.03%  "; ForceText = $false }
    @{ Cell = "E38"; Value = "  -2.44%  "; ForceText = $false }
    @{ Cell = "E39"; Value = "  -1.04%  "; ForceText = $false }
    @{ Cell = "D40"; Value = "0.101"; ForceText = $true }
    @{ Cell = "E40"; Value = "  -3.66%  "; ForceText = $false }
    @{ Cell = "E41"; Value = "  +0.24%  "; ForceText = $false }
    @{ Cell = "E42"; Value = "  -0.35%  "; ForceText = $false }
    @{ Cell = "D43"; Value = "19.57"; ForceText = $true }
    @{ Cell = "E43"; Value = "  -8.25%  "; ForceText = $false }
    @{ Cell = "D44"; Value = "1.950.26"; ForceText = $false }
    @{ Cell = "E44"; Value = "  -1.59%  "; ForceText = $false }
    @{ Cell = "E45"; Value = "  -0.51%  "; ForceText = $false }
    @{ Cell = "D46"; Value = "2.97"; ForceText = $true }
    @{ Cell = "E46"; Value = "  -2.79%  "; ForceText = $false }
    @{ Cell = "D47"; Value = "8.80"; ForceText = $true }
    @{ Cell = "E47"; Value = "  +0.24%  "; ForceText = $false }
    @{ Cell = "D48"; Value = "2.727.54"; ForceText = $false }
    @{ Cell = "E48"; Value = "  +0.96%  "; ForceText = $false }
    @{ Cell = "D49"; Value = "96.04"; ForceText = $true }
    @{ Cell = "E49"; Value = "  -1.50%  "; ForceText = $false }
    @{ Cell = "D50"; Value = "0.176"; ForceText = $true }
    @{ Cell = "E50"; Value = "  -2.10%  "; ForceText = $false }
    @{ Cell = "D51"; Value = "67.15"; ForceText = $true }
    @{ Cell = "E51"; Value = "  -2.65%  "; ForceText = $false }
)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    if ($u.ForceText) {
        # These values look numeric (e.g. "0.996", "96.04") but the source
        # data is plain text, so force a text format before assigning the
        # value to stop Excel from auto-converting/rounding it to a number,
        # then drop back to the default style so no formatting is changed.
        $rng.NumberFormat = "@"
        $rng.Value = $u.Value
        $rng.Style = "Normal"
    } else {
        $rng.Value = $u.Value
    }
}
